$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 43786
$ws.Range("E2").Value = 856663314701
$ws.Range("F2").Value = 13429325067
$ws.Range("G2").Value = -0.05802

$ws.Range("D3").Value = 2342.8
$ws.Range("E3").Value = 281679393489
$ws.Range("F3").Value = 11447335822
$ws.Range("G3").Value = -0.53816

$ws.Range("E4").Value = 90634565590
$ws.Range("F4").Value = 32628239434
$ws.Range("G4").Value = 0.0265

$ws.Range("D5").Value = 239.15
$ws.Range("E5").Value = 36771511247
$ws.Range("F5").Value = 663402131
$ws.Range("G5").Value = -1.2681

$ws.Range("D6").Value = 0.659102
$ws.Range("E6").Value = 35544704216
$ws.Range("F6").Value = 1390707472
$ws.Range("G6").Value = -3.71253

$ws.Range("D7").Value = 72.17
$ws.Range("E7").Value = 30721567472
$ws.Range("F7").Value = 2248614651
$ws.Range("G7").Value = -6.31586

$ws.Range("E8").Value = 24551289562
$ws.Range("F8").Value = 6588430696
$ws.Range("G8").Value = -0.08277

$ws.Range("D9").Value = 2339.66
$ws.Range("E9").Value = 21581083249
$ws.Range("F9").Value = 24010751
$ws.Range("G9").Value = -0.42705

$ws.Range("D10").Value = 0.596977
$ws.Range("E10").Value = 20821338858
$ws.Range("F10").Value = 1544977249
$ws.Range("G10").Value = -5.00953

$ws.Range("D11").Value = 0.099148
$ws.Range("E11").Value = 14084536479
$ws.Range("F11").Value = 1147059410
$ws.Range("G11").Value = -2.57604

$ws.Range("B12").Value = "AVAX"
$ws.Range("C12").Value = "Avalanche"
$ws.Range("D12").Value = 32.56
$ws.Range("E12").Value = 11867937817
$ws.Range("F12").Value = 1202243194
$ws.Range("G12").Value = -3.29568

$ws.Range("B13").Value = "TRX"
$ws.Range("C13").Value = "TRON"
$ws.Range("D13").Value = 0.107163
$ws.Range("E13").Value = 9483396758
$ws.Range("F13").Value = 342331969
$ws.Range("G13").Value = -0.56198

$ws.Range("B14").Value = "DOT"
$ws.Range("C14").Value = "Polkadot"
$ws.Range("D14").Value = 7.2
$ws.Range("E14").Value = 9371740976
$ws.Range("F14").Value = 385992652
$ws.Range("G14").Value = -5.29052

$ws.Range("B15").Value = "LINK"
$ws.Range("C15").Value = "Chainlink"
$ws.Range("D15").Value = 16.03
$ws.Range("E15").Value = 8908056779
$ws.Range("F15").Value = 711143092
$ws.Range("G15").Value = -5.41228

$ws.Range("D16").Value = 0.898182
$ws.Range("E16").Value = 8338188820
$ws.Range("F16").Value = 610106446
$ws.Range("G16").Value = -3.12373

$ws.Range("B17").Value = "TON"
$ws.Range("C17").Value = "Toncoin"
$ws.Range("D17").Value = 2.31
$ws.Range("E17").Value = 8121723900
$ws.Range("F17").Value = 11666911
$ws.Range("G17").Value = -0.64437

$ws.Range("D18").Value = 43741
$ws.Range("E18").Value = 6766780585
$ws.Range("F18").Value = 149131786
$ws.Range("G18").Value = 0.15431

$ws.Range("B19").Value = "SHIB"
$ws.Range("C19").Value = "Shiba Inu"
$ws.Range("D19").Value = 0.00001007
$ws.Range("E19").Value = 5932161616
$ws.Range("F19").Value = 446837436
$ws.Range("G19").Value = -1.77626

$ws.Range("B20").Value = "LTC"
$ws.Range("C20").Value = "Litecoin"
$ws.Range("D20").Value = 77.70999999999999
$ws.Range("E20").Value = 5750231658
$ws.Range("F20").Value = 782209266
$ws.Range("G20").Value = 0.08064

$ws.Range("B21").Value = "DAI"
$ws.Range("C21").Value = "Dai"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 5423264512
$ws.Range("F21").Value = 249835221
$ws.Range("G21").Value = 0.00762

$ws.Range("D22").Value = 6.62
$ws.Range("E22").Value = 4987631838
$ws.Range("F22").Value = 232568980
$ws.Range("G22").Value = -0.98506

$ws.Range("D23").Value = 248.9
$ws.Range("E23").Value = 4869103156
$ws.Range("F23").Value = 149665351
$ws.Range("G23").Value = -2.48626

$ws.Range("B24").Value = "XLM"
$ws.Range("C24").Value = "Stellar"
$ws.Range("D24").Value = 0.132588
$ws.Range("E24").Value = 3735848768
$ws.Range("F24").Value = 97186014
$ws.Range("G24").Value = -1.80773

$ws.Range("D25").Value = 3.79
$ws.Range("E25").Value = 3520894905
$ws.Range("F25").Value = 1143483
$ws.Range("G25").Value = -0.16781

$ws.Range("B26").Value = "OKB"
$ws.Range("C26").Value = "OKB"
$ws.Range("D26").Value = 58.12
$ws.Range("E26").Value = 3484801297
$ws.Range("F26").Value = 46371021
$ws.Range("G26").Value = 1.2512

$ws.Range("D27").Value = 176.89
$ws.Range("E27").Value = 3202302250
$ws.Range("F27").Value = 68883053
$ws.Range("G27").Value = 0.7456700000000001

$ws.Range("B28").Value = "ETC"
$ws.Range("C28").Value = "Ethereum Classic"
$ws.Range("D28").Value = 21.88
$ws.Range("E28").Value = 3130812756
$ws.Range("F28").Value = 197032243
$ws.Range("G28").Value = -3.20925

$ws.Range("B29").Value = "ATOM"
$ws.Range("C29").Value = "Cosmos Hub"
$ws.Range("D29").Value = 10.29
$ws.Range("E29").Value = 3007229065
$ws.Range("F29").Value = 175846988
$ws.Range("G29").Value = -6.37044

$ws.Range("B30").Value = "CRO"
$ws.Range("C30").Value = "Cronos"
$ws.Range("D30").Value = 0.105902
$ws.Range("E30").Value = 2801076291
$ws.Range("F30").Value = 51937612
$ws.Range("G30").Value = -0.19628

$ws.Range("B31").Value = "KAS"
$ws.Range("C31").Value = "Kaspa"
$ws.Range("D31").Value = 0.125818
$ws.Range("E31").Value = 2740235633
$ws.Range("F31").Value = 30193224
$ws.Range("G31").Value = -2.3427

$ws.Range("B32").Value = "TUSD"
$ws.Range("C32").Value = "TrueUSD"
$ws.Range("D32").Value = 0.998962
$ws.Range("E32").Value = 2630853649
$ws.Range("F32").Value = 126764272
$ws.Range("G32").Value = 0.26312

$ws.Range("B33").Value = "HBAR"
$ws.Range("C33").Value = "Hedera"
$ws.Range("D33").Value = 0.073156
$ws.Range("E33").Value = 2454658936
$ws.Range("F33").Value = 63871760
$ws.Range("G33").Value = -4.00019

$ws.Range("B34").Value = "FIL"
$ws.Range("C34").Value = "Filecoin"
$ws.Range("D34").Value = 5.09
$ws.Range("E34").Value = 2428199572
$ws.Range("F34").Value = 224424030
$ws.Range("G34").Value = -5.33712

$ws.Range("B35").Value = "NEAR"
$ws.Range("C35").Value = "NEAR Protocol"
$ws.Range("D35").Value = 2.41
$ws.Range("E35").Value = 2420473466
$ws.Range("F35").Value = 250522557
$ws.Range("G35").Value = -4.86508

$ws.Range("B36").Value = "ICP"
$ws.Range("C36").Value = "Internet Computer"
$ws.Range("D36").Value = 5.34
$ws.Range("E36").Value = 2402137200
$ws.Range("F36").Value = 67952558
$ws.Range("G36").Value = -1.12428

$ws.Range("B37").Value = "IMX"
$ws.Range("C37").Value = "Immutable"
$ws.Range("D37").Value = 1.83
$ws.Range("E37").Value = 2351447661
$ws.Range("F37").Value = 1308624556
$ws.Range("G37").Value = 2.90395

$ws.Range("D38").Value = 8.18
$ws.Range("E38").Value = 2284748587
$ws.Range("F38").Value = 145919568
$ws.Range("G38").Value = -4.86736

$ws.Range("B39").Value = "LDO"
$ws.Range("C39").Value = "Lido DAO"
$ws.Range("D39").Value = 2.37
$ws.Range("E39").Value = 2103184514
$ws.Range("F39").Value = 55662851
$ws.Range("G39").Value = -2.88263

$ws.Range("B40").Value = "OP"
$ws.Range("C40").Value = "Optimism"
$ws.Range("D40").Value = 2.2
$ws.Range("E40").Value = 1992984289
$ws.Range("F40").Value = 291256819
$ws.Range("G40").Value = 4.25876

$ws.Range("B41").Value = "VET"
$ws.Range("C41").Value = "VeChain"
$ws.Range("D41").Value = 0.02687952
$ws.Range("E41").Value = 1952221059
$ws.Range("F41").Value = 59769243
$ws.Range("G41").Value = -2.54601

$ws.Range("B42").Value = "TAO"
$ws.Range("C42").Value = "Bittensor"
$ws.Range("D42").Value = 335.69
$ws.Range("E42").Value = 1950759039
$ws.Range("F42").Value = 4056577
$ws.Range("G42").Value = -1.0043

$ws.Range("B43").Value = "RUNE"
$ws.Range("C43").Value = "THORChain"
$ws.Range("D43").Value = 6.39
$ws.Range("E43").Value = 1918702746
$ws.Range("F43").Value = 416539843
$ws.Range("G43").Value = -1.1742

$ws.Range("B44").Value = "MNT"
$ws.Range("C44").Value = "Mantle"
$ws.Range("D44").Value = 0.606142
$ws.Range("E44").Value = 1899717987
$ws.Range("F44").Value = 6424115
$ws.Range("G44").Value = -2.71496

$ws.Range("B45").Value = "EGLD"
$ws.Range("C45").Value = "MultiversX"
$ws.Range("D45").Value = 64.55
$ws.Range("E45").Value = 1702130997
$ws.Range("F45").Value = 163415276
$ws.Range("G45").Value = 15.21462

$ws.Range("B46").Value = "QNT"
$ws.Range("C46").Value = "Quant"
$ws.Range("D46").Value = 115.47
$ws.Range("E46").Value = 1677477962
$ws.Range("F46").Value = 30136318
$ws.Range("G46").Value = -2.6346

$ws.Range("D47").Value = 18.84
$ws.Range("E47").Value = 1581668185
$ws.Range("F47").Value = 94653951
$ws.Range("G47").Value = -2.69951

$ws.Range("B48").Value = "GRT"
$ws.Range("C48").Value = "The Graph"
$ws.Range("D48").Value = 0.168841
$ws.Range("E48").Value = 1572126142
$ws.Range("F48").Value = 74881659
$ws.Range("G48").Value = -1.58061

$ws.Range("B49").Value = "ALGO"
$ws.Range("C49").Value = "Algorand"
$ws.Range("D49").Value = 0.193841
$ws.Range("E49").Value = 1546298373
$ws.Range("F49").Value = 212675441
$ws.Range("G49").Value = -4.17292

$ws.Range("B50").Value = "BUSD"
$ws.Range("C50").Value = "BUSD"
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 1506606703
$ws.Range("F50").Value = 2040555489
$ws.Range("G50").Value = 0.03979

$ws.Range("B51").Value = "ARB"
$ws.Range("C51").Value = "Arbitrum"
$ws.Range("D51").Value = 1.14
$ws.Range("E51").Value = 1448577421
$ws.Range("F51").Value = 375115468
$ws.Range("G51").Value = -4.03384
